$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Row 1 (header): relabel the existing columns and append the new schema columns ---
$ws1.Range("B1").Value = "name"
$ws1.Range("C1").Value = "capacity"
$ws1.Range("D1").Value = "owner"
$ws1.Range("E1").Value = "register_date"
$ws1.Range("F1").Value = "register_reason"
$ws1.Range("G1").Value = "acquire_value"
$ws1.Range("H1").Value = "property_category"
$ws1.Range("I1").Value = "category"
$ws1.Range("J1").Value = "date"
$ws1.Range("K1").Value = "legislator_name"
$ws1.Range("L1").Value = "legislator_id"
$ws1.Range("M1").Value = "source_file"
$ws1.Range("N1").Value = "index"

# Match the bold/bordered header formatting (same style as B1:G1) on the new header cells
$ws1.Range("G1").Copy() | Out-Null
$ws1.Range("H1:N1").PasteSpecial(-4122) | Out-Null

# --- Row 2 (data): append values for the new schema columns ---
$ws1.Range("H2").Value = "land"
$ws1.Range("I2").Value = "normal"
# Force text so "2012-05-01" is kept as a literal string, not auto-converted to a date serial
$ws1.Range("J2").NumberFormat = "@"
$ws1.Range("J2").Value = "2012-05-01"
$ws1.Range("K2").Value = "陳亭妃"
$ws1.Range("L2").Value = 1708
$ws1.Range("M2").Value = "tmpb2a21"
$ws1.Range("N2").Value = 29

# Match the plain/bordered data-row formatting (same style as B2:G2) on the new data cells
$ws1.Range("G2").Copy() | Out-Null
$ws1.Range("H2:N2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
